$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix import functionality: re-populate the R-xxx identifier column (A2:A22)
# with freshly (re)generated values, now emitted in sorted R-number order
# instead of the prior scrambled order. (Mirrors "newAltIndexToAdd" rebuild.)
$ws.Range("A2").Value = "R-201-I1-CS1-7B4D"
$ws.Range("A3").Value = "R-202-I2-Cb1-F317"
$ws.Range("A4").Value = "R-203-I3-CS2-5B61"
$ws.Range("A5").Value = "R-204-I4-Cf1-62A0"
$ws.Range("A6").Value = "R-205-I5-Ce1-4C99"
$ws.Range("A7").Value = "R-206-I6-CT1-64E3"
$ws.Range("A8").Value = "R-207-I7-CE2-9616"
$ws.Range("A9").Value = "R-208-I8-CS1-D002"
$ws.Range("A10").Value = "R-209-I9-Cb1-0D0C"
$ws.Range("A11").Value = "R-210-I10-CS2-8217"
$ws.Range("A12").Value = "R-211-I11-Cf1-ACC5"
$ws.Range("A13").Value = "R-212-I12-Ce1-00F4"
$ws.Range("A14").Value = "R-213-I13-CT1-B5AB"
$ws.Range("A15").Value = "R-214-I14-CE2-1A71"
$ws.Range("A16").Value = "R-215-I15-CS1-4111"
$ws.Range("A17").Value = "R-216-I16-Cb1-7510"
$ws.Range("A18").Value = "R-217-I17-CS2-942F"
$ws.Range("A19").Value = "R-218-I18-Cf1-AE28"
$ws.Range("A20").Value = "R-219-I19-Ce1-8FFF"
$ws.Range("A21").Value = "R-220-I20-CT1-6919"
$ws.Range("A22").Value = "R-222-I22-CS2-06EC"

# Leave the editor's cursor where the author's session ended up.
$ws.Range("A23").Select()
